# Apply the "Updated cryptos list" refresh: new prices / 1h-volume deltas for
# every coin row, plus a handful of rows whose rank order changed (their B/C/D/E
# cells got swapped with the neighbouring row).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '91.592.25'
$ws.Range('E2').Value = '  +2.10%  '

# Row 3
$ws.Range('D3').Value = '3.105.44'
$ws.Range('E3').Value = '  +1.00%  '

# Row 4
$ws.Range('E4').Value = '  +0.02%  '

# Row 5
$ws.Range('D5').Value = '''242.72'
$ws.Range('E5').Value = '  +1.41%  '

# Row 6
$ws.Range('D6').Value = '''616.00'
$ws.Range('E6').Value = '  -0.13%  '

# Row 7
$ws.Range('E7').Value = '  -3.93%  '

# Row 8
$ws.Range('D8').Value = '''0.393'
$ws.Range('E8').Value = '  +8.69%  '

# Row 9
$ws.Range('E9').Value = '  -0.04%  '

# Row 10
$ws.Range('D10').Value = '3.105.93'
$ws.Range('E10').Value = '  +1.04%  '

# Row 11
$ws.Range('E11').Value = '  +0.33%  '

# Row 12
$ws.Range('E12').Value = '  +0.24%  '

# Row 13
$ws.Range('D13').Value = '''0.0000252'
$ws.Range('E13').Value = '  +3.61%  '

# Row 14
$ws.Range('D14').Value = '''34.52'
$ws.Range('E14').Value = '  +0.39%  '

# Row 15
$ws.Range('D15').Value = '92.048.49'
$ws.Range('E15').Value = '  +2.63%  '

# Row 16
$ws.Range('D16').Value = '''5.52'
$ws.Range('E16').Value = '  +1.40%  '

# Row 17
$ws.Range('D17').Value = '3.689.75'
$ws.Range('E17').Value = '  +0.97%  '

# Row 18
$ws.Range('D18').Value = '3.158.04'
$ws.Range('E18').Value = '  +2.30%  '

# Row 19
$ws.Range('D19').Value = '''3.64'
$ws.Range('E19').Value = '  -0.26%  '

# Row 20
$ws.Range('E20').Value = '  +2.63%  '

# Row 21
$ws.Range('D21').Value = '''5.81'
$ws.Range('E21').Value = '  +1.54%  '

# Row 22
$ws.Range('D22').Value = '''447.27'
$ws.Range('E22').Value = '  +2.85%  '

# Row 23 (Uniswap)
$ws.Range('B23').Value = 'Uniswap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D23').Value = '''9.29'
$ws.Range('E23').Value = '  +3.92%  '

# Row 24 (PEPE)
$ws.Range('B24').Value = 'PEPE'
$ws.Range('C24').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D24').Value = '''0.0000202'
$ws.Range('E24').Value = '  -2.13%  '

# Row 25
$ws.Range('D25').Value = '''5.63'

# Row 26 (Litecoin)
$ws.Range('B26').Value = 'Litecoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D26').Value = '''87.01'
$ws.Range('E26').Value = '  -3.51%  '

# Row 27 (Aptos)
$ws.Range('B27').Value = 'Aptos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D27').Value = '''11.65'
$ws.Range('E27').Value = '  -0.37%  '

# Row 28
$ws.Range('D28').Value = '3.278.69'

# Row 29
$ws.Range('E29').Value = '  -0.34%  '

# Row 30
$ws.Range('D30').Value = '''0.139'
$ws.Range('E30').Value = '  +20.55%  '

# Row 31
$ws.Range('E31').Value = '  -5.52%  '

# Row 32
$ws.Range('D32').Value = '''0.167'
$ws.Range('E32').Value = '  -4.97%  '

# Row 33
$ws.Range('E33').Value = '  +2.84%  '

# Row 34 (Binance-PegBSC-USD)
$ws.Range('B34').Value = 'Binance-PegBSC-USD'
$ws.Range('C34').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D34').Value = '''1.00'
$ws.Range('E34').Value = '  +3.40%  '

# Row 35 (Kaspa)
$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').Value = '''0.169'
$ws.Range('E35').Value = '  +2.79%  '

# Row 36
$ws.Range('E36').Value = '  +4.85%  '

# Row 37
$ws.Range('D37').Value = '''26.25'
$ws.Range('E37').Value = '  +0.80%  '

# Row 38
$ws.Range('D38').Value = '''4.10'
$ws.Range('E38').Value = '  -3.97%  '

# Row 39
$ws.Range('D39').Value = '''1.94'
$ws.Range('E39').Value = '  +2.16%  '

# Row 40
$ws.Range('E40').Value = '  +2.79%  '

# Row 41
$ws.Range('D41').Value = '''480.39'
$ws.Range('E41').Value = '  -0.31%  '

# Row 42 (dogwifhat)
$ws.Range('B42').Value = 'dogwifhat'
$ws.Range('C42').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D42').Value = '''3.46'
$ws.Range('E42').Value = '  -0.60%  '

# Row 43 (PolygonEcosystemToken)
$ws.Range('B43').Value = 'PolygonEcosystemToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D43').Value = '''0.433'
$ws.Range('E43').Value = '  +4.60%  '

# Row 44
$ws.Range('D44').Value = '''22.19'
$ws.Range('E44').Value = '  +0.27%  '

# Row 45
$ws.Range('E45').Value = '  -0.01%  '

# Row 46
$ws.Range('E46').Value = '  +2.80%  '

# Row 47
$ws.Range('D47').Value = '''1.91'
$ws.Range('E47').Value = '  +2.22%  '

# Row 48
$ws.Range('E48').Value = '  +2.94%  '

# Row 49
$ws.Range('E49').Value = '  +2.80%  '

# Row 50
$ws.Range('D50').Value = '''0.0336'
$ws.Range('E50').Value = '  +9.42%  '

# Row 51
$ws.Range('D51').Value = '''43.92'
$ws.Range('E51').Value = '  -0.21%  '

# The quote-prefix assignments above mark those cells as "Text" and stamp a
# quotePrefix style on them; restore the default "Normal" style so formatting
# matches the rest of the untouched column (value/text-ness is unaffected).
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').Style = 'Normal'
